$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns (D:E, data rows) to Text format first so that
# numeric-looking strings (e.g. "249.80") are NOT auto-converted to numbers,
# preserving the original inline-string cell type.
$ws.Range('D2:E51').NumberFormat = "@"

$ws.Range('D2').Value2 = '96.770.22'
$ws.Range('E2').Value2 = '  +2.65%  '
$ws.Range('D3').Value2 = '3.301.74'
$ws.Range('E3').Value2 = '  +6.24%  '
$ws.Range('E4').Value2 = '  -0.15%  '
$ws.Range('D5').Value2 = '249.80'
$ws.Range('E5').Value2 = '  +6.01%  '
$ws.Range('D6').Value2 = '619.31'
$ws.Range('E6').Value2 = '  +1.21%  '
$ws.Range('D7').Value2 = '1.10'
$ws.Range('E7').Value2 = '  -0.43%  '
$ws.Range('D8').Value2 = '0.381'
$ws.Range('E8').Value2 = '  -1.04%  '
$ws.Range('E9').Value2 = '  -0.05%  '
$ws.Range('D10').Value2 = '3.300.03'
$ws.Range('E10').Value2 = '  +6.30%  '
$ws.Range('D11').Value2 = '0.784'
$ws.Range('E11').Value2 = '  -5.73%  '
$ws.Range('E12').Value2 = '  +0.57%  '
$ws.Range('D13').Value2 = '96.223.14'
$ws.Range('E13').Value2 = '  +2.24%  '
$ws.Range('D14').Value2 = '0.0000245'
$ws.Range('E14').Value2 = '  +0.95%  '
$ws.Range('D15').Value2 = '35.13'
$ws.Range('E15').Value2 = '  +2.04%  '
$ws.Range('D16').Value2 = '3.875.11'
$ws.Range('E16').Value2 = '  +4.89%  '
$ws.Range('D17').Value2 = '5.50'
$ws.Range('E17').Value2 = '  +5.17%  '
$ws.Range('D18').Value2 = '3.279.83'
$ws.Range('E18').Value2 = '  +5.16%  '
$ws.Range('D19').Value2 = '3.56'
$ws.Range('E19').Value2 = '  -2.78%  '
$ws.Range('D20').Value2 = '14.88'
$ws.Range('E20').Value2 = '  +0.92%  '
$ws.Range('D21').Value2 = '477.39'
$ws.Range('E21').Value2 = '  +6.98%  '
$ws.Range('D22').Value2 = '5.80'
$ws.Range('E22').Value2 = '  -0.74%  '
$ws.Range('D23').Value2 = '0.0000206'
$ws.Range('E23').Value2 = '  +4.86%  '
$ws.Range('D24').Value2 = '9.24'
$ws.Range('E24').Value2 = '  +3.49%  '
$ws.Range('D25').Value2 = '5.63'
$ws.Range('E25').Value2 = '  +0.88%  '
$ws.Range('D26').Value2 = '87.52'
$ws.Range('E26').Value2 = '  +2.17%  '
$ws.Range('D27').Value2 = '12.17'
$ws.Range('E27').Value2 = '  +0.02%  '
$ws.Range('D28').Value2 = '3.454.90'
$ws.Range('E28').Value2 = '  +5.01%  '
$ws.Range('E29').Value2 = '  +0.14%  '
$ws.Range('D30').Value2 = '0.182'
$ws.Range('E30').Value2 = '  -2.14%  '
$ws.Range('D31').Value2 = '0.240'
$ws.Range('E31').Value2 = '  -7.42%  '
$ws.Range('D32').Value2 = '0.999'
$ws.Range('E32').Value2 = '  +0.96%  '
$ws.Range('D33').Value2 = '0.121'
$ws.Range('E33').Value2 = '  -1.05%  '
$ws.Range('D34').Value2 = '9.21'
$ws.Range('E34').Value2 = '  -0.74%  '
$ws.Range('D35').Value2 = '27.03'
$ws.Range('E35').Value2 = '  +4.66%  '
$ws.Range('D36').Value2 = '7.40'
$ws.Range('E36').Value2 = '  -6.15%  '
$ws.Range('E37').Value2 = '  -5.57%  '
$ws.Range('B38').Value2 = 'Bittensor'
$ws.Range('C38').Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value2 = '500.87'
$ws.Range('E38').Value2 = '  +6.79%  '
$ws.Range('B39').Value2 = 'PancakeSwap'
$ws.Range('C39').Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D39').Value2 = '1.93'
$ws.Range('E39').Value2 = '  +1.67%  '
$ws.Range('B40').Value2 = 'WhiteBITCoin'
$ws.Range('C40').Value2 = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D40').Value2 = '24.73'
$ws.Range('E40').Value2 = '  +3.28%  '
$ws.Range('D41').Value2 = '0.447'
$ws.Range('E41').Value2 = '  -0.67%  '
$ws.Range('D42').Value2 = '1.27'
$ws.Range('E42').Value2 = '  -0.71%  '
$ws.Range('D43').Value2 = '3.26'
$ws.Range('E43').Value2 = '  +2.87%  '
$ws.Range('D44').Value2 = '0.788'
$ws.Range('E44').Value2 = '  +14.77%  '
$ws.Range('B45').Value2 = 'MantraDAO'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D45').Value2 = '3.48'
$ws.Range('E45').Value2 = '  -4.32%  '
$ws.Range('B46').Value2 = 'USDe'
$ws.Range('C46').Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value2 = '1.00'
$ws.Range('E46').Value2 = '  -0.03%  '
$ws.Range('D47').Value2 = '160.97'
$ws.Range('E47').Value2 = '  -0.27%  '
$ws.Range('D48').Value2 = '1.90'
$ws.Range('E48').Value2 = '  +1.74%  '
$ws.Range('B49').Value2 = 'ImmutableX'
$ws.Range('C49').Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').Value2 = '1.37'
$ws.Range('E49').Value2 = '  +5.24%  '
$ws.Range('B50').Value2 = 'OKB'
$ws.Range('C50').Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value2 = '45.22'
$ws.Range('E50').Value2 = '  +3.41%  '
$ws.Range('D51').Value2 = '4.49'
$ws.Range('E51').Value2 = '  +3.09%  '

# Restore default styling on the range (removes the temporary text format,
# cell values remain text since they were already committed as strings).
$ws.Range('D2:E51').Style = "Normal"

